# Update the latest automated-test-run results on the active sheet
# (Plan1) of the scenarios backup workbook: refresh the run timestamps
# and flip a couple of outcome cells to reflect the newest test pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (CT 01): refresh run timestamp
$ws.Range("H2").Value = "28_04_2020--21_27_33 556"

# Row 5 (CT 04): outcome flips from Yes to No, timestamp refreshed
$ws.Range("B5").Value = "No"
$ws.Range("H5").Value = "28_04_2020--21_22_15 376"

# Row 6 (CT 05): status flips from Passed to Failed, timestamp refreshed
$ws.Range("C6").Value = "Failed"
$ws.Range("H6").Value = "28_04_2020--21_28_27 839"

# Move the active selection to B10 (beyond the current data, matching
# where the next scenario row will be appended)
$ws.Range("B10").Select()
